# Regenerate merged AHB files
# 1) Rename header labels from *_old / *_new to *_FV2304 / *_FV2310
# 2) Turn the data range into a native Excel Table ("Table1")
# 3) Freeze the header row (pane split after row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..J carry the "_old" -> "_FV2304" headers
$colsOld = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $colsOld) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2304")
}

# Columns L..U carry the "_new" -> "_FV2310" headers (K is "diff", unchanged)
$colsNew = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $colsNew) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2310")
}

# Convert the used range A1:U94 into a native Excel table ("Table1")
$tableRange = $ws.Range("A1:U94")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row (split/freeze after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
